$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.198.42"
$ws.Range("E2").Value = "  +8.57%  "
$ws.Range("D3").Value = "2.428.51"
$ws.Range("E3").Value = "  +9.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "478.49"
$ws.Range("E5").Value = "  +13.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.76"
$ws.Range("E6").Value = "  +24.91%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("E8").Value = "  +14.39%  "
$ws.Range("D9").Value = "2.452.40"
$ws.Range("E9").Value = "  +10.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0959"
$ws.Range("E10").Value = "  +18.73%  "
$ws.Range("E11").Value = "  +8.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.322"
$ws.Range("E12").Value = "  +12.80%  "
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").Value = "2.842.51"
$ws.Range("E14").Value = "  +8.57%  "
$ws.Range("D15").Value = "55.150.48"
$ws.Range("E15").Value = "  +8.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.46"
$ws.Range("E16").Value = "  +16.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000134"
$ws.Range("E17").Value = "  +22.97%  "
$ws.Range("D18").Value = "2.445.86"
$ws.Range("E18").Value = "  +8.60%  "
$ws.Range("E19").Value = "  +16.39%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.83"
$ws.Range("E20").Value = "  +20.45%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.10"
$ws.Range("E21").Value = "  +10.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.993"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("E23").Value = "  +19.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.06"
$ws.Range("E24").Value = "  +10.84%  "
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.400"
$ws.Range("E26").Value = "  +15.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").Value = "  +22.59%  "
$ws.Range("D28").Value = "2.542.85"
$ws.Range("E28").Value = "  +9.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.33"
$ws.Range("E29").Value = "  +14.75%  "
$ws.Range("D30").Value = "0.0₃0766"
$ws.Range("E30").Value = "  +30.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.48"
$ws.Range("E32").Value = "  +4.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.87"
$ws.Range("E33").Value = "  +12.56%  "
$ws.Range("E34").Value = "  +16.55%  "
$ws.Range("E35").Value = "  +16.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.12"
$ws.Range("E36").Value = "  +20.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.57"
$ws.Range("E37").Value = "  +13.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.840"
$ws.Range("E38").Value = "  +12.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.43"
$ws.Range("E39").Value = "  +7.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.994"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.601"
$ws.Range("E41").Value = "  +11.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  +15.98%  "
$ws.Range("E43").Value = "  +15.26%  "
$ws.Range("E44").Value = "  +20.00%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "254.23"
$ws.Range("E46").Value = "  +39.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.61"
$ws.Range("E47").Value = "  +25.10%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0222"
$ws.Range("E48").Value = "  +15.90%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0890"
$ws.Range("E49").Value = "  +15.68%  "
$ws.Range("D50").Value = "1.914.46"
$ws.Range("E50").Value = "  +5.39%  "
$ws.Range("E51").Value = "  +14.98%  "
